$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.813.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.74%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.859.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = '  +0.86%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("E5").Value = '  +0.72%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = '  +0.74%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4413"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.98%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3817"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.76%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.07445"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.60%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.8882"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.43%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'21.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.871.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.99%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'5.540"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'6.749"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.07217"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.55%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = '  +4.27%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'1.041"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.96%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.000009114"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.71%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'1.035"
$ws.Range("D19").Style = "Normal"

$ws.Range("E20").Value = '  +0.94%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'27.821.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.77%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'5.308"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.94%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'11.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.40%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'2.089.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.68%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = '  +6.97%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'159.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.18%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'18.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.24%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'1.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.10%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'5.359"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.97%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'118.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.26%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.09119"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.58%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.7734"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.52%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = '  +0.31%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'3.030"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.40%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'4.610"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.87%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'1.158"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.28%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = '  +0.40%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.05320"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.62%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'2.865"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.79%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.5217"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.64%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'6.960"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.41%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = '  +0.34%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'8.790"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.08%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'110.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.69%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'10.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.39%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = '  +0.87%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = '  +2.99%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.4733"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.53%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'1.879"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.16%  '
$ws.Range("E51").Style = "Normal"
